# Recompute NATMI LR-pair metrics (Sema3g -> Nrp2) with updated TPM inputs.
# Ligand-expressing-cell counts / detection rates / derived specificity scores
# all shift for rows 2-10 (ECs, FAPs, MuSCs sender x ECs/FAPs/MuSCs target).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.538602
$ws.Range("H2").Value = 22.615806
$ws.Range("I2").Value = 0.8744677208212333
$ws.Range("J2").Value = 0.8744677208212334
$ws.Range("M2").Value = 18.43631966666667
$ws.Range("N2").Value = 55.308959
$ws.Range("O2").Value = 0.6034704469962782
$ws.Range("P2").Value = 0.6034704469962781
$ws.Range("Q2").Value = 138.9840763117727
$ws.Range("R2").Value = 1250.856686805954
$ws.Range("S2").Value = 0.5277154263678062
$ws.Range("T2").Value = 0.5277154263678062

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.538602
$ws.Range("H3").Value = 22.615806
$ws.Range("I3").Value = 0.8744677208212333
$ws.Range("J3").Value = 0.8744677208212334
$ws.Range("O3").Value = 0.1750419652256785
$ws.Range("P3").Value = 0.1750419652256784
$ws.Range("Q3").Value = 40.31356626290333
$ws.Range("R3").Value = 362.82209636613
$ws.Range("S3").Value = 0.1530685483789686
$ws.Range("T3").Value = 0.1530685483789686

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.538602
$ws.Range("H4").Value = 22.615806
$ws.Range("I4").Value = 0.8744677208212333
$ws.Range("J4").Value = 0.8744677208212334
$ws.Range("M4").Value = 6.766555
$ws.Range("N4").Value = 20.299665
$ws.Range("O4").Value = 0.2214875877780434
$ws.Range("P4").Value = 0.2214875877780434
$ws.Range("Q4").Value = 51.01036505611
$ws.Range("R4").Value = 459.09328550499
$ws.Range("S4").Value = 0.1936837460744585
$ws.Range("T4").Value = 0.1936837460744585

# Row 5
$ws.Range("I5").Value = 0.03298905189192265
$ws.Range("J5").Value = 0.03298905189192266
$ws.Range("M5").Value = 18.43631966666667
$ws.Range("N5").Value = 55.308959
$ws.Range("O5").Value = 0.6034704469962782
$ws.Range("P5").Value = 0.6034704469962781
$ws.Range("Q5").Value = 5.243135677202777
$ws.Range("R5").Value = 47.188221094825
$ws.Range("S5").Value = 0.01990791789120198
$ws.Range("T5").Value = 0.01990791789120198

# Row 6
$ws.Range("I6").Value = 0.03298905189192265
$ws.Range("J6").Value = 0.03298905189192266
$ws.Range("O6").Value = 0.1750419652256785
$ws.Range("P6").Value = 0.1750419652256784
$ws.Range("S6").Value = 0.005774468474094028
$ws.Range("T6").Value = 0.005774468474094028

# Row 7
$ws.Range("I7").Value = 0.03298905189192265
$ws.Range("J7").Value = 0.03298905189192266
$ws.Range("M7").Value = 6.766555
$ws.Range("N7").Value = 20.299665
$ws.Range("O7").Value = 0.2214875877780434
$ws.Range("P7").Value = 0.2214875877780434
$ws.Range("Q7").Value = 1.924351854041667
$ws.Range("R7").Value = 17.319166686375
$ws.Range("S7").Value = 0.007306665526626647
$ws.Range("T7").Value = 0.007306665526626647

# Row 8
$ws.Range("G8").Value = 0.7977956666666666
$ws.Range("H8").Value = 2.393387
$ws.Range("I8").Value = 0.09254322728684393
$ws.Range("J8").Value = 0.09254322728684394
$ws.Range("M8").Value = 18.43631966666667
$ws.Range("N8").Value = 55.308959
$ws.Range("O8").Value = 0.6034704469962782
$ws.Range("P8").Value = 0.6034704469962781
$ws.Range("Q8").Value = 14.70841593934811
$ws.Range("R8").Value = 132.375743454133
$ws.Range("S8").Value = 0.05584710273726987
$ws.Range("T8").Value = 0.05584710273726987

# Row 9
$ws.Range("G9").Value = 0.7977956666666666
$ws.Range("H9").Value = 2.393387
$ws.Range("I9").Value = 0.09254322728684393
$ws.Range("J9").Value = 0.09254322728684394
$ws.Range("O9").Value = 0.1750419652256785
$ws.Range("P9").Value = 0.1750419652256784
$ws.Range("Q9").Value = 4.266306733320555
$ws.Range("R9").Value = 38.39676059988499
$ws.Range("S9").Value = 0.01619894837261579
$ws.Range("T9").Value = 0.01619894837261579

# Row 10
$ws.Range("G10").Value = 0.7977956666666666
$ws.Range("H10").Value = 2.393387
$ws.Range("I10").Value = 0.09254322728684393
$ws.Range("J10").Value = 0.09254322728684394
$ws.Range("M10").Value = 6.766555
$ws.Range("N10").Value = 20.299665
$ws.Range("O10").Value = 0.2214875877780434
$ws.Range("P10").Value = 0.2214875877780434
$ws.Range("Q10").Value = 5.398328257261666
$ws.Range("R10").Value = 48.584954315355
$ws.Range("S10").Value = 0.02049717617695827
$ws.Range("T10").Value = 0.02049717617695827
